$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete rows 8-10 (table shrank from 9 data rows to 6 data rows)
$ws.Range("A8:T10").EntireRow.Delete()

# Update remaining data rows (2-7) with refreshed TPM-derived values
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Anxa1"
$ws.Range("C2").Value = "Fpr1"
$ws.Range("D2").Value = "FAPs"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 2.811979666666667
$ws.Range("H2").Value = 8.435939000000001
$ws.Range("I2").Value = 0.01221198172659148
$ws.Range("J2").Value = 0.01221198172659148
$ws.Range("K2").Value = 1
$ws.Range("L2").Value = 0.3333333333333333
$ws.Range("M2").Value = 0.6032380000000001
$ws.Range("N2").Value = 1.809714
$ws.Range("O2").Value = 0.8880288532312676
$ws.Range("P2").Value = 0.8880288532312675
$ws.Range("Q2").Value = 1.696292990160667
$ws.Range("R2").Value = 15.266636911446
$ws.Range("S2").Value = 0.01084459212834623
$ws.Range("T2").Value = 0.01084459212834623
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Anxa1"
$ws.Range("C3").Value = "Fpr1"
$ws.Range("D3").Value = "MuSCs"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 2.811979666666667
$ws.Range("H3").Value = 8.435939000000001
$ws.Range("I3").Value = 0.01221198172659148
$ws.Range("J3").Value = 0.01221198172659148
$ws.Range("K3").Value = 1
$ws.Range("L3").Value = 0.3333333333333333
$ws.Range("M3").Value = 0.076062
$ws.Range("N3").Value = 0.228186
$ws.Range("O3").Value = 0.1119711467687325
$ws.Range("P3").Value = 0.1119711467687325
$ws.Range("Q3").Value = 0.2138847974060001
$ws.Range("R3").Value = 1.924963176654
$ws.Range("S3").Value = 0.001367389598245254
$ws.Range("T3").Value = 0.001367389598245254
$ws.Range("A4").Value = "FAPs"
$ws.Range("B4").Value = "Anxa1"
$ws.Range("C4").Value = "Fpr1"
$ws.Range("D4").Value = "FAPs"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 186.4134573333333
$ws.Range("H4").Value = 559.240372
$ws.Range("I4").Value = 0.8095640809678946
$ws.Range("J4").Value = 0.8095640809678947
$ws.Range("K4").Value = 1
$ws.Range("L4").Value = 0.3333333333333333
$ws.Range("M4").Value = 0.6032380000000001
$ws.Range("N4").Value = 1.809714
$ws.Range("O4").Value = 0.8880288532312676
$ws.Range("P4").Value = 0.8880288532312675
$ws.Range("Q4").Value = 112.4516811748453
$ws.Range("R4").Value = 1012.065130573608
$ws.Range("S4").Value = 0.7189162624391444
$ws.Range("T4").Value = 0.7189162624391444
$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "Anxa1"
$ws.Range("C5").Value = "Fpr1"
$ws.Range("D5").Value = "MuSCs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 186.4134573333333
$ws.Range("H5").Value = 559.240372
$ws.Range("I5").Value = 0.8095640809678946
$ws.Range("J5").Value = 0.8095640809678947
$ws.Range("K5").Value = 1
$ws.Range("L5").Value = 0.3333333333333333
$ws.Range("M5").Value = 0.076062
$ws.Range("N5").Value = 0.228186
$ws.Range("O5").Value = 0.1119711467687325
$ws.Range("P5").Value = 0.1119711467687325
$ws.Range("Q5").Value = 14.178980391688
$ws.Range("R5").Value = 127.610823525192
$ws.Range("S5").Value = 0.09064781852875019
$ws.Range("T5").Value = 0.09064781852875019
$ws.Range("A6").Value = "MuSCs"
$ws.Range("B6").Value = "Anxa1"
$ws.Range("C6").Value = "Fpr1"
$ws.Range("D6").Value = "FAPs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 41.03855533333333
$ws.Range("H6").Value = 123.115666
$ws.Range("I6").Value = 0.1782239373055139
$ws.Range("J6").Value = 0.1782239373055139
$ws.Range("K6").Value = 1
$ws.Range("L6").Value = 0.3333333333333333
$ws.Range("M6").Value = 0.6032380000000001
$ws.Range("N6").Value = 1.809714
$ws.Range("O6").Value = 0.8880288532312676
$ws.Range("P6").Value = 0.8880288532312675
$ws.Range("Q6").Value = 24.75601604216934
$ws.Range("R6").Value = 222.804144379524
$ws.Range("S6").Value = 0.1582679986637769
$ws.Range("T6").Value = 0.1582679986637768
$ws.Range("A7").Value = "MuSCs"
$ws.Range("B7").Value = "Anxa1"
$ws.Range("C7").Value = "Fpr1"
$ws.Range("D7").Value = "MuSCs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 41.03855533333333
$ws.Range("H7").Value = 123.115666
$ws.Range("I7").Value = 0.1782239373055139
$ws.Range("J7").Value = 0.1782239373055139
$ws.Range("K7").Value = 1
$ws.Range("L7").Value = 0.3333333333333333
$ws.Range("M7").Value = 0.076062
$ws.Range("N7").Value = 0.228186
$ws.Range("O7").Value = 0.1119711467687325
$ws.Range("P7").Value = 0.1119711467687325
$ws.Range("Q7").Value = 3.121474595764
$ws.Range("R7").Value = 28.093271361876
$ws.Range("S7").Value = 0.01995593864173708
$ws.Range("T7").Value = 0.01995593864173708
